$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Donor")

# Insert a new column at Y, shifting the existing "LoadedAt" column (Y) to Z
$ws.Columns("Y:Y").Insert()

# New column header
$ws.Range("Y1").Value = "AccumulatedDonationValue"

# New column values (AccumulatedDonationValue)
$ws.Range("Y2").Value = 800
$ws.Range("Y3").Value = 9000
$ws.Range("Y4").Value = 324
$ws.Range("Y5").Value = 56.32
$ws.Range("Y6").Value = 123.45

# Match target column width (OOXML width="22")
$ws.Columns("Y:Y").ColumnWidth = 21.17

# Update the view: scroll so column Q is the left-most visible column
$excel.ActiveWindow.ScrollColumn = 17
$excel.ActiveWindow.ScrollRow = 1

# Select Y7 to match the saved selection state
$ws.Range("Y7").Select()
